$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for three new rows (Megadrive, Genesis - Multitap, Megadrive - Multitap)
# right after the existing "Genesis" row (row 15), by shifting rows 16-30 down to
# rows 19-33. We move row-by-row from the bottom up (Cut/Paste) instead of using
# Rows.Insert() so that the moved cells keep their original style indexes without
# Excel fabricating a brand-new (unused) cell format in styles.xml.
for ($r = 30; $r -ge 16; $r--) {
    $srcRow = $r
    $dstRow = $r + 3
    $src = $ws.Range("A$srcRow" + ":D$srcRow")
    $dst = $ws.Range("A$dstRow" + ":D$dstRow")
    $src.Cut($dst)
}

# Fill in the three newly freed rows with the Multitap entries. The shared-string
# table records new unique strings in first-use order, so write the cells in the
# same order the authored workbook introduces them: "Genesis - Multitap",
# "Megadrive - Multitap", then "Megadrive".
$ws.Range("A17").Value = "Genesis - Multitap"
$ws.Range("A18").Value = "Megadrive - Multitap"
$ws.Range("A16").Value = "Megadrive"

$ws.Range("B16").Value = "Complete"
$ws.Range("C16").Value = "Complete"

$ws.Range("B17").Value = "Complete"
$ws.Range("C17").Value = "Complete"

$ws.Range("B18").Value = "Complete"
$ws.Range("C18").Value = "Complete"

# Update the active selection to match the authored workbook.
$ws.Range("J25").Select() | Out-Null
